$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText "D2" "70.740.60"
Set-CellText "E2" "  +1.93%  "
Set-CellText "D3" "3.470.83"
Set-CellText "E3" "  +2.35%  "
Set-CellText "E4" "  -0.05%  "
Set-CellText "D5" "587.65"
Set-CellText "E5" "  +0.03%  "
Set-CellText "D6" "179.58"
Set-CellText "E6" "  -0.14%  "
Set-CellText "B7" "LidoStakedEther"
Set-CellText "C7" "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-CellText "D7" "3.462.63"
Set-CellText "E7" "  +2.29%  "
Set-CellText "B8" "XRP"
Set-CellText "C8" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-CellText "D8" "0.603"
Set-CellText "E8" "  +1.13%  "
Set-CellText "E10" "  +6.06%  "
Set-CellText "D11" "0.594"
Set-CellText "E11" "  +0.53%  "
Set-CellText "D12" "49.43"
Set-CellText "E12" "  +1.90%  "
Set-CellText "D13" "0.0000288"
Set-CellText "E13" "  +1.98%  "
Set-CellText "D14" "695.09"
Set-CellText "E14" "  +2.44%  "
Set-CellText "B15" "Polkadot"
Set-CellText "C15" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText "D15" "8.77"
Set-CellText "E15" "  +1.80%  "
Set-CellText "B16" "WrappedliquidstakedEther2.0"
Set-CellText "C16" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText "D16" "4.019.82"
Set-CellText "E16" "  +2.08%  "
Set-CellText "D17" "70.622.86"
Set-CellText "E17" "  +1.68%  "
Set-CellText "D18" "3.447.16"
Set-CellText "E18" "  +1.61%  "
Set-CellText "D19" "0.122"
Set-CellText "E19" "  +1.06%  "
Set-CellText "D20" "17.88"
Set-CellText "E20" "  +1.38%  "
Set-CellText "D21" "11.53"
Set-CellText "E21" "  +2.36%  "
Set-CellText "D22" "0.912"
Set-CellText "E22" "  +1.04%  "
Set-CellText "D23" "5.51"
Set-CellText "E23" "  +1.65%  "
Set-CellText "D24" "17.19"
Set-CellText "E24" "  +0.53%  "
Set-CellText "D25" "101.83"
Set-CellText "E25" "  -1.55%  "
Set-CellText "D26" "3.97"
Set-CellText "E26" "  +1.09%  "
Set-CellText "D27" "2.71"
Set-CellText "E27" "  -0.57%  "
Set-CellText "D28" "9.74"
Set-CellText "E28" "  +0.65%  "
Set-CellText "D29" "33.98"
Set-CellText "E29" "  -0.48%  "
Set-CellText "D30" "8.89"
Set-CellText "E30" "  +2.09%  "
Set-CellText "D31" "7.24"
Set-CellText "E31" "  +3.69%  "
Set-CellText "D32" "3.94"
Set-CellText "E32" "  +9.03%  "
Set-CellText "D33" "577.58"
Set-CellText "E33" "  +3.92%  "
Set-CellText "D34" "11.12"
Set-CellText "E34" "  -0.33%  "
Set-CellText "D35" "59.04"
Set-CellText "E35" "  +1.74%  "
Set-CellText "E36" "  -2.39%  "
Set-CellText "E37" "  +0.09%  "
Set-CellText "D38" "3.603.36"
Set-CellText "E38" "  -2.24%  "
Set-CellText "D39" "0.141"
Set-CellText "E39" "  +1.06%  "
Set-CellText "D40" "35.62"
Set-CellText "E40" "  +1.75%  "
Set-CellText "D41" "0.0₃0749"
Set-CellText "E41" "  +7.25%  "
Set-CellText "D42" "3.37"
Set-CellText "E42" "  +2.90%  "
Set-CellText "D43" "2.74"
Set-CellText "E43" "  +2.09%  "
Set-CellText "B44" "TheGraph"
Set-CellText "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-CellText "D44" "0.340"
Set-CellText "E44" "  +0.41%  "
Set-CellText "D45" "0.0427"
Set-CellText "E45" "  +1.06%  "
Set-CellText "B46" "ApeXProtocol"
Set-CellText "C46" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-CellText "D46" "3.35"
Set-CellText "E46" "  +1.41%  "
Set-CellText "D47" "2.72"
Set-CellText "E47" "  +1.67%  "
Set-CellText "E48" "  +3.43%  "
Set-CellText "E49" "  +0.53%  "
Set-CellText "D50" "0.996"
Set-CellText "E50" "  -0.48%  "
Set-CellText "D51" "133.73"
Set-CellText "E51" "  +0.82%  "
